$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the scraped cryptos list refresh
$ws.Range("D2").Value = '26.269.56'
$ws.Range("E2").Value = '  -2.91%  '
$ws.Range("D3").Value = '1.770.37'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.98'
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4218'
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3593'
$ws.Range("E8").Value = '  +1.32%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07116'
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8353'
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.30'
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").Value = '1.762.23'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.431'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.225'
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06878'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '78.59'
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008663'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Value = '26.285.68'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.088'
$ws.Range("E22").Value = '  +0.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.94'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '1.983.17'
$ws.Range("E24").Value = '  -1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.81'
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.791'
$ws.Range("E26").Value = '  -7.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.95'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.043'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.89'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.832'
$ws.Range("E30").Value = '  +11.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08826'
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7226'
$ws.Range("E32").Value = '  +1.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.110'
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.291'
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("E35").Value = '  -0.31%  '
$ws.Range("E36").Value = '  -4.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.110'
$ws.Range("E37").Value = '  +3.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05093'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01880'
$ws.Range("E39").Value = '  -0.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.1604'
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4895'
$ws.Range("E41").Value = '  -0.53%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.590'
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.318'
$ws.Range("E43").Value = '  +3.95%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.998'
$ws.Range("E44").Value = '  -0.03%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.25'
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.45'
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.613'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06160'
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4448'
$ws.Range("E50").Value = '  -0.91%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.686'
$ws.Range("E51").Value = '  +1.11%  '
